# Remove errant trailing semicolon from "_h2o_keep_element;" text runs
# (found in paragraphs styled "NodeEnd"), leaving plain "_h2o_keep_element".

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Execute a Find & Replace All for the exact errant text.
$find.Execute(
    "_h2o_keep_element;",  # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "_h2o_keep_element",   # ReplaceWith
    2                      # Replace (wdReplaceAll)
)
